$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.675.40"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.278.47"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.74"
$ws.Range("E5").Value = "  +3.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.94"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.274.98"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.175"
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.571"
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.22"
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").Value = "3.805.55"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "620.29"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").Value = "65.562.71"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.80"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "3.280.84"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.90"
$ws.Range("E21").Value = "  -3.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.888"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.05"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.77"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.97"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.52"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.85"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.43"
$ws.Range("E30").Value = "  -2.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.42"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.72"
$ws.Range("E32").Value = "  -8.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "546.98"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.84"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("D35").Value = "3.789.39"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.03"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.45"
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("E41").Value = "  +3.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.59"
$ws.Range("E43").Value = "  -4.09%  "
$ws.Range("D44").Value = "0.0₃0677"
$ws.Range("E44").Value = "  -8.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.330"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0405"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.00"
$ws.Range("E47").Value = "  -6.56%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.50"
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.29"
$ws.Range("E51").Value = "  +5.09%  "
